$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '61.691.53'
$ws.Range("E2").Value = '  +1.53%  '
$ws.Range("D3").Value = '3.399.68'
$ws.Range("E3").Value = '  +1.10%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.67'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.80%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '3.399.76'
$ws.Range("E8").Value = '  +1.12%  '
$ws.Range("E9").Value = '  -0.78%  '
$ws.Range("E10").Value = '  -0.64%  '
$ws.Range("E11").Value = '  +3.03%  '
$ws.Range("E12").Value = '  +0.02%  '
$ws.Range("D13").Value = '3.976.87'
$ws.Range("E13").Value = '  +1.01%  '
$ws.Range("E14").Value = '  +2.18%  '
$ws.Range("E15").Value = '  +2.92%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.94'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.12%  '
$ws.Range("D17").Value = '3.397.00'
$ws.Range("E17").Value = '  +1.00%  '
$ws.Range("D18").Value = '61.822.59'
$ws.Range("E18").Value = '  +1.36%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.24'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.47%  '
$ws.Range("E20").Value = '  +1.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.47'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '378.07'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.60%  '
$ws.Range("E23").Value = '  -1.35%  '
$ws.Range("D24").Value = '3.528.77'
$ws.Range("E24").Value = '  +0.93%  '
$ws.Range("B25").Value = 'Dai'
$ws.Range("C25").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("B26").Value = 'PEPE'
$ws.Range("C26").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000127'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +9.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '71.28'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.97%  '
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.57'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.55%  '
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.161'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.66%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.26'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.88%  '
$ws.Range("E33").Value = '  +1.97%  '
$ws.Range("E34").Value = '  +0.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.46'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.35%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.35'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.73%  '
$ws.Range("E37").Value = '  +1.07%  '
$ws.Range("E38").Value = '  -0.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '165.38'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.45%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0783'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("E41").Value = '  +9.16%  '
$ws.Range("B42").Value = 'Mantle'
$ws.Range("C42").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.782'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.14%  '
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.07%  '
$ws.Range("B44").Value = 'ONDO'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.23'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.66%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '25.16'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +9.38%  '
$ws.Range("E46").Value = '  +0.74%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '41.50'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.86'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.79'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.71%  '
$ws.Range("D50").Value = '2.342.33'
$ws.Range("E50").Value = '  +5.99%  '
$ws.Range("E51").Value = '  +2.52%  '
